$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Arrived" column (H) as "Y" (arrived) for the items that have now
# come in.
$ws.Range("H3").Value = "Y"
$ws.Range("H6").Value = "Y"
$ws.Range("H31").Value = "Y"
$ws.Range("H35").Value = "Y"
$ws.Range("H38").Value = "Y"
$ws.Range("H42").Value = "Y"
$ws.Range("H43").Value = "Y"

# Update the view: scroll back to the top and move the active selection to H4.
$ws.Range("A1").Select()
$ws.Range("H4").Select()
